$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume/coin data
# Force text format on cells so numeric-looking strings (e.g. prices)
# are preserved exactly as text rather than being parsed as numbers.
$cell = $ws.Cells.Item(2, 4)
$cell.NumberFormat = "@"
$cell.Value = '63.872.75'
$cell = $ws.Cells.Item(2, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -2.79%  '
$cell = $ws.Cells.Item(3, 4)
$cell.NumberFormat = "@"
$cell.Value = '3.167.18'
$cell = $ws.Cells.Item(3, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -3.28%  '
$cell = $ws.Cells.Item(4, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.05%  '
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = '566.52'
$cell = $ws.Cells.Item(5, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -2.86%  '
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = '167.67'
$cell = $ws.Cells.Item(6, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -6.54%  '
$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.605'
$cell = $ws.Cells.Item(7, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -6.86%  '
$cell = $ws.Cells.Item(8, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +0.00%  '
$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = '3.167.68'
$cell = $ws.Cells.Item(9, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -3.35%  '
$cell = $ws.Cells.Item(10, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -3.44%  '
$cell = $ws.Cells.Item(11, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.77%  '
$cell = $ws.Cells.Item(12, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -4.72%  '
$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = '3.719.49'
$cell = $ws.Cells.Item(13, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -3.32%  '
$cell = $ws.Cells.Item(14, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -1.81%  '
$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = '63.971.11'
$cell = $ws.Cells.Item(15, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -2.83%  '
$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = '25.28'
$cell = $ws.Cells.Item(16, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -2.80%  '
$cell = $ws.Cells.Item(17, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -2.72%  '
$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = '3.177.20'
$cell = $ws.Cells.Item(18, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -2.81%  '
$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = '413.31'
$cell = $ws.Cells.Item(19, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -3.19%  '
$cell = $ws.Cells.Item(20, 2)
$cell.NumberFormat = "@"
$cell.Value = 'Chainlink'
$cell = $ws.Cells.Item(20, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = '12.75'
$cell = $ws.Cells.Item(20, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -3.33%  '
$cell = $ws.Cells.Item(21, 2)
$cell.NumberFormat = "@"
$cell.Value = 'Polkadot'
$cell = $ws.Cells.Item(21, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.32'
$cell = $ws.Cells.Item(21, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -2.86%  '
$cell = $ws.Cells.Item(22, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -3.78%  '
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell = $ws.Cells.Item(23, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +0.01%  '
$cell = $ws.Cells.Item(24, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -1.82%  '
$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.202'
$cell = $ws.Cells.Item(25, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +3.20%  '
$cell = $ws.Cells.Item(26, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -3.73%  '
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0000107'
$cell = $ws.Cells.Item(27, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -5.27%  '
$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = "@"
$cell.Value = '8.68'
$cell = $ws.Cells.Item(28, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -2.07%  '
$cell = $ws.Cells.Item(29, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.13%  '
$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.83'
$cell = $ws.Cells.Item(30, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -6.85%  '
$cell = $ws.Cells.Item(31, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -2.14%  '
$cell = $ws.Cells.Item(32, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.03%  '
$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = '4.99'
$cell = $ws.Cells.Item(33, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -3.05%  '
$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = '6.33'
$cell = $ws.Cells.Item(34, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -3.66%  '
$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.12'
$cell = $ws.Cells.Item(35, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -4.90%  '
$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = "@"
$cell.Value = '155.84'
$cell = $ws.Cells.Item(36, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -2.07%  '
$cell = $ws.Cells.Item(37, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -3.37%  '
$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.730.21'
$cell = $ws.Cells.Item(38, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -1.73%  '
$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.69'
$cell = $ws.Cells.Item(39, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -5.40%  '
$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = '24.64'
$cell = $ws.Cells.Item(40, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -6.33%  '
$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = '4.15'
$cell = $ws.Cells.Item(41, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -3.80%  '
$cell = $ws.Cells.Item(42, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -7.11%  '
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = '38.63'
$cell = $ws.Cells.Item(43, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -3.53%  '
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0621'
$cell = $ws.Cells.Item(44, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -5.42%  '
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.59'
$cell = $ws.Cells.Item(45, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -5.46%  '
$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0261'
$cell = $ws.Cells.Item(46, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -2.04%  '
$cell = $ws.Cells.Item(47, 2)
$cell.NumberFormat = "@"
$cell.Value = 'InjectiveProtocol'
$cell = $ws.Cells.Item(47, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = '21.76'
$cell = $ws.Cells.Item(47, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -5.52%  '
$cell = $ws.Cells.Item(48, 2)
$cell.NumberFormat = "@"
$cell.Value = 'Bittensor'
$cell = $ws.Cells.Item(48, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = '293.64'
$cell = $ws.Cells.Item(48, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -6.59%  '
$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.04'
$cell = $ws.Cells.Item(50, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.10%  '
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0987'
$cell = $ws.Cells.Item(51, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -6.34%  '
